{"js": "// Update the visit date (Date: field) and the \"Next review\" date.\n// \"05 / Oct / 2021\" -> \"06 / Oct / 2021\"\n// \"19 / 10 / 2021\"  -> \"20 / 10 / 2021\"\n\nconst body = context.document.body;\n\nconst dateResults = body.search(\"05 / Oct / 2021\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items/text\");\nawait context.sync();\n\nfor (const range of dateResults.items) {\n  range.insertText(\"06 / Oct / 2021\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst reviewResults = body.search(\"19 / 10 / 2021\", { matchCase: true, matchWholeWord: false });\nreviewResults.load(\"items/text\");\nawait context.sync();\n\nfor (const range of reviewResults.items) {\n  range.insertText(\"20 / 10 / 2021\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the visit date (Date: field) and the \"Next review\" date.\n# \"05 / Oct / 2021\" -> \"06 / Oct / 2021\"\n# \"19 / 10 / 2021\"  -> \"20 / 10 / 2021\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"05 / Oct / 2021\"\n$find.Replacement.Text = \"06 / Oct / 2021\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"19 / 10 / 2021\"\n$find2.Replacement.Text = \"20 / 10 / 2021\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
